$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits at the end of the paragraph
#    that ends in "...properly." (just before its paragraph mark). The
#    edit moves that bookmark to the end of the very last paragraph in
#    the document (the "Implement Butler-Reed-Dawson..." paragraph).
#    Remove it from its old spot first.
# -----------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# -----------------------------------------------------------------------
# 2. Fix the misspelling "vesion" -> "version" in the last paragraph.
#    The target markup drops the spell-check proofErr wrapper entirely
#    and ends up with the corrected word split across three separate
#    runs: "...simplified ve" | "r" | "sion of the method...".
#    Delete the whole paragraph's text (this clears out the proofErr
#    elements) and retype the corrected sentence.
# -----------------------------------------------------------------------
$lastIdx = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIdx)

$bodyRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$bodyRange.Delete()

$lastPara = $d.Paragraphs.Item($lastIdx)
$insPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insPoint.InsertAfter("Implement Butler-Reed-Dawson algorithm and a simplified version of the method in the paper  with 1 dimension for T2")

# Re-split the freshly typed text into three runs at the "ve|r" and
# "r|sion" boundaries. Dropping a (temporary) bookmark at a position
# forces the engine to break the run there; removing the bookmark
# again afterwards leaves the break in place without adding any
# leftover run formatting.
$lastPara = $d.Paragraphs.Item($lastIdx)
$paraStart = $lastPara.Range.Start

$splitAfterVe = $d.Range($paraStart + 58, $paraStart + 58)
$d.Bookmarks.Add("_TmpRunSplit1", $splitAfterVe)

$splitAfterR = $d.Range($paraStart + 59, $paraStart + 59)
$d.Bookmarks.Add("_TmpRunSplit2", $splitAfterR)

$d.Bookmarks("_TmpRunSplit1").Delete()
$d.Bookmarks("_TmpRunSplit2").Delete()

# -----------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark at the end of this (now last)
#    paragraph, right before its paragraph mark - matching where it
#    ended up in the target document.
# -----------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($lastIdx)
$newBookmarkRange = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
